$d = $word.ActiveDocument

# Mapping of old text -> new text, applied via Find/Replace across the
# whole document. Each old value is unique in the document, so a simple
# Replace-All for each pair is safe and order-independent.
$replacements = @(
    @("2026-01-15 Thursday", "2026-01-16 Friday"),
    @("742÷8=", "370÷8="),
    @("241÷2=", "464÷7="),
    @("590÷2=", "489÷3="),
    @("343÷3=", "914÷6="),
    @("606÷6=", "865÷7="),
    @("182÷6=", "274÷4="),
    @("628÷2=", "685÷4="),
    @("609÷5=", "384÷6="),
    @("888÷7=", "207÷2="),
    @("975÷9=", "816÷7="),
    @("454÷9=", "605÷5="),
    @("515÷8=", "853÷6="),
    @("417÷6=", "217÷9="),
    @("397÷6=", "662÷8="),
    @("992÷8=", "828÷7="),
    @("102÷4=", "483÷4="),
    @("143÷4=", "107÷3="),
    @("809÷4=", "455÷5="),
    @("798÷8=", "968÷2="),
    @("312÷5=", "528÷8="),
    @("864÷4=", "394÷8="),
    @("135÷9=", "230÷7="),
    @("508÷3=", "804÷2="),
    @("649÷4=", "154÷7="),
    @("312÷2=", "402÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
